$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.744.96'
$ws.Range('E2').Value = '  +0.89%  '
$ws.Range('D3').Value = '2.317.63'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'319.96"
$ws.Range('E5').Value = '  +2.43%  '
$ws.Range('D6').Value = "'104.27"
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').Value = "'40.07"
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('D11').Value = "'0.0913"
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').Value = "'8.38"
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').Value = "'0.976"
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = "'15.45"
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('D16').Value = '2.666.53'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '2.325.81'
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('D18').Value = '42.676.55'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').Value = "'7.54"
$ws.Range('E19').Value = '  +0.42%  '
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').Value = "'3.64"
$ws.Range('E21').Value = '  +3.46%  '
$ws.Range('D22').Value = "'73.57"
$ws.Range('E22').Value = '  -2.57%  '
$ws.Range('D23').Value = "'279.31"
$ws.Range('E23').Value = '  +6.48%  '
$ws.Range('D24').Value = "'11.10"
$ws.Range('E24').Value = '  +19.19%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('D27').Value = "'10.94"
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('D28').Value = "'2.37"
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('D29').Value = "'23.01"
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').Value = "'36.22"
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').Value = "'165.31"
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = "'0.0882"
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('D33').Value = "'5.97"
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('E34').Value = '  +5.73%  '
$ws.Range('E35').Value = '  -8.93%  '
$ws.Range('E36').Value = '  -3.37%  '
$ws.Range('D37').Value = "'4.66"
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('D38').Value = "'0.0360"
$ws.Range('E38').Value = '  +2.24%  '
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('E40').Value = '  +4.72%  '
$ws.Range('E41').Value = '  +2.53%  '
$ws.Range('D42').Value = "'99.30"
$ws.Range('E42').Value = '  -1.35%  '
$ws.Range('D43').Value = "'70.07"
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('E44').Value = '  -2.11%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'114.17"
$ws.Range('E47').Value = '  +2.15%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = "'80.70"
$ws.Range('E48').Value = '  +7.97%  '
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('D51').Value = '1.615.76'
$ws.Range('E51').Value = '  +4.67%  '
